# Apply targeted cell-value updates across all 8 job sheets (Moogle_Profits workbook).
# Values correspond to recalculated market-price-driven profit figures.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1937.9166
$ws.Range("I5").Value = 1937.9166
$ws.Range("K5").Value = 1937.9166
$ws.Range("M5").Value = -1822.9166
$ws.Range("H11").Value = 131.5625
$ws.Range("I11").Value = 131.5625
$ws.Range("K11").Value = 131.5625
$ws.Range("M11").Value = 8.4375
$ws.Range("H38").Value = 1420.3478
$ws.Range("I38").Value = 209.33333
$ws.Range("K38").Value = 627.99999
$ws.Range("M38").Value = -255.99999
$ws.Range("H86").Value = 4849.1934
$ws.Range("I86").Value = 2701.5264
$ws.Range("J86").Value = 8249.667
$ws.Range("K86").Value = 2701.5264
$ws.Range("L86").Value = 8249.667
$ws.Range("M86").Value = -1578.5264
$ws.Range("N86").Value = -10495.667
$ws.Range("H89").Value = 4849.1934
$ws.Range("I89").Value = 2701.5264
$ws.Range("J89").Value = 8249.667
$ws.Range("K89").Value = 13507.632
$ws.Range("L89").Value = 41248.335
$ws.Range("M89").Value = -7891.632000000001
$ws.Range("N89").Value = -52480.335
$ws.Range("H98").Value = 3067.375
$ws.Range("I98").Value = 2077
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 2077
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -579
$ws.Range("N98").Value = -12996
$ws.Range("H104").Value = 469.66666
$ws.Range("I104").Value = 163.6
$ws.Range("K104").Value = 490.8
$ws.Range("M104").Value = 1256.2
$ws.Range("H106").Value = 20954198
$ws.Range("I106").Value = 24445914
$ws.Range("K106").Value = 24445914
$ws.Range("M106").Value = -24445283
$ws.Range("H112").Value = 4170.5713
$ws.Range("J112").Value = 4170.5713
$ws.Range("L112").Value = 12511.7139
$ws.Range("N112").Value = -14727.7139
$ws.Range("H122").Value = 3067.375
$ws.Range("I122").Value = 2077
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 6231
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -3781
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 2807.3965
$ws.Range("I132").Value = 2559.2263
$ws.Range("K132").Value = 7677.678899999999
$ws.Range("M132").Value = -5147.678899999999
$ws.Range("H137").Value = 6312.6865
$ws.Range("I137").Value = 5640.3057
$ws.Range("K137").Value = 16920.9171
$ws.Range("M137").Value = -14370.9171
$ws.Range("H138").Value = 5821.6665
$ws.Range("I138").Value = 4358.091
$ws.Range("J138").Value = 6224.15
$ws.Range("K138").Value = 13074.273
$ws.Range("L138").Value = 18672.45
$ws.Range("M138").Value = -7934.273000000001
$ws.Range("N138").Value = -28952.45
$ws.Range("H141").Value = 3943.5925
$ws.Range("I141").Value = 2415.6667
$ws.Range("K141").Value = 7247.000100000001
$ws.Range("M141").Value = -2067.000100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5287.9
$ws.Range("I32").Value = 4191.2954
$ws.Range("K32").Value = 4191.2954
$ws.Range("M32").Value = -3904.2954
$ws.Range("H41").Value = 4999.5
$ws.Range("I41").Value = 4999.5
$ws.Range("K41").Value = 4999.5
$ws.Range("M41").Value = -4585.5
$ws.Range("H74").Value = 4518.423
$ws.Range("I74").Value = 2422.8096
$ws.Range("J74").Value = 13320
$ws.Range("K74").Value = 2422.8096
$ws.Range("L74").Value = 13320
$ws.Range("M74").Value = -1548.8096
$ws.Range("N74").Value = -15068
$ws.Range("H77").Value = 4518.423
$ws.Range("I77").Value = 2422.8096
$ws.Range("J77").Value = 13320
$ws.Range("K77").Value = 12114.048
$ws.Range("L77").Value = 66600
$ws.Range("M77").Value = -7746.048000000001
$ws.Range("N77").Value = -75336
$ws.Range("H80").Value = 105950
$ws.Range("J80").Value = 122437.5
$ws.Range("L80").Value = 122437.5
$ws.Range("N80").Value = -124433.5
$ws.Range("H83").Value = 105950
$ws.Range("J83").Value = 122437.5
$ws.Range("L83").Value = 367312.5
$ws.Range("N83").Value = -377296.5
$ws.Range("H97").Value = 1591.9375
$ws.Range("I97").Value = 1668.4667
$ws.Range("J97").Value = 444
$ws.Range("K97").Value = 1668.4667
$ws.Range("L97").Value = 444
$ws.Range("M97").Value = -1172.4667
$ws.Range("N97").Value = -1436
$ws.Range("H122").Value = 2392.182
$ws.Range("I122").Value = 2468.8572
$ws.Range("K122").Value = 7406.571599999999
$ws.Range("M122").Value = -4956.571599999999
$ws.Range("H132").Value = 3250.7966
$ws.Range("I132").Value = 2192.1372
$ws.Range("J132").Value = 9999.75
$ws.Range("K132").Value = 6576.4116
$ws.Range("L132").Value = 29999.25
$ws.Range("M132").Value = -4046.4116
$ws.Range("N132").Value = -35059.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 78875
$ws.Range("J35").Value = 78875
$ws.Range("L35").Value = 78875
$ws.Range("N35").Value = -79495
$ws.Range("H82").Value = 47718.625
$ws.Range("I82").Value = 12795
$ws.Range("J82").Value = 152489.5
$ws.Range("K82").Value = 12795
$ws.Range("L82").Value = 152489.5
$ws.Range("M82").Value = -12412
$ws.Range("N82").Value = -153255.5
$ws.Range("H85").Value = 47718.625
$ws.Range("I85").Value = 12795
$ws.Range("J85").Value = 152489.5
$ws.Range("K85").Value = 12795
$ws.Range("L85").Value = 152489.5
$ws.Range("M85").Value = -11469
$ws.Range("N85").Value = -155141.5
$ws.Range("H86").Value = 6479
$ws.Range("I86").Value = 3418.9
$ws.Range("K86").Value = 3418.9
$ws.Range("M86").Value = -2295.9
$ws.Range("H89").Value = 6479
$ws.Range("I89").Value = 3418.9
$ws.Range("K89").Value = 17094.5
$ws.Range("M89").Value = -11478.5
$ws.Range("H94").Value = 1393
$ws.Range("I94").Value = 948.2857
$ws.Range("K94").Value = 948.2857
$ws.Range("M94").Value = -497.2857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3600
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -913
$ws.Range("H25").Value = 5693.154
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 5000
$ws.Range("N25").Value = -5348
$ws.Range("H31").Value = 8161.027
$ws.Range("I31").Value = 4606.6333
$ws.Range("K31").Value = 4606.6333
$ws.Range("M31").Value = -4311.6333
$ws.Range("H34").Value = 8161.027
$ws.Range("I34").Value = 4606.6333
$ws.Range("K34").Value = 4606.6333
$ws.Range("M34").Value = -4404.6333
$ws.Range("H41").Value = 16399.6
$ws.Range("J41").Value = 19332.666
$ws.Range("L41").Value = 19332.666
$ws.Range("N41").Value = -20188.666
$ws.Range("H50").Value = 51875
$ws.Range("J50").Value = 51875
$ws.Range("L50").Value = 51875
$ws.Range("N50").Value = -53125
$ws.Range("H58").Value = 7460.55
$ws.Range("I58").Value = 4815
$ws.Range("J58").Value = 13633.5
$ws.Range("K58").Value = 4815
$ws.Range("L58").Value = 13633.5
$ws.Range("M58").Value = -4612
$ws.Range("N58").Value = -14039.5
$ws.Range("H59").Value = 99165.78
$ws.Range("I59").Value = 25000
$ws.Range("J59").Value = 108436.5
$ws.Range("K59").Value = 25000
$ws.Range("L59").Value = 108436.5
$ws.Range("M59").Value = -23855
$ws.Range("N59").Value = -110726.5
$ws.Range("H60").Value = 27281.834
$ws.Range("J60").Value = 31538.2
$ws.Range("L60").Value = 31538.2
$ws.Range("N60").Value = -32560.2
$ws.Range("H86").Value = 7226.222
$ws.Range("I86").Value = 6506.1665
$ws.Range("J86").Value = 8666.333
$ws.Range("K86").Value = 6506.1665
$ws.Range("L86").Value = 8666.333
$ws.Range("M86").Value = -5383.1665
$ws.Range("N86").Value = -10912.333
$ws.Range("H89").Value = 7226.222
$ws.Range("I89").Value = 6506.1665
$ws.Range("J89").Value = 8666.333
$ws.Range("K89").Value = 32530.8325
$ws.Range("L89").Value = 43331.665
$ws.Range("M89").Value = -26914.8325
$ws.Range("N89").Value = -54563.665
$ws.Range("H94").Value = 3709.818
$ws.Range("I94").Value = 2878.2
$ws.Range("K94").Value = 2878.2
$ws.Range("M94").Value = -2427.2
$ws.Range("H113").Value = 3600
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970
$ws.Range("H122").Value = 1281.4615
$ws.Range("I122").Value = 1031.8096
$ws.Range("K122").Value = 3095.4288
$ws.Range("M122").Value = -645.4288000000001
$ws.Range("H125").Value = 44999.5
$ws.Range("J125").Value = 44999.5
$ws.Range("L125").Value = 44999.5
$ws.Range("N125").Value = -49919.5
$ws.Range("H132").Value = 5589.393
$ws.Range("I132").Value = 3840.5293
$ws.Range("K132").Value = 11521.5879
$ws.Range("M132").Value = -8991.5879
$ws.Range("H134").Value = 3963.68
$ws.Range("I134").Value = 2636.5264
$ws.Range("J134").Value = 8166.3335
$ws.Range("K134").Value = 7909.5792
$ws.Range("L134").Value = 24499.0005
$ws.Range("M134").Value = -5374.5792
$ws.Range("N134").Value = -29569.0005
$ws.Range("H136").Value = 7460.55
$ws.Range("I136").Value = 4815
$ws.Range("J136").Value = 13633.5
$ws.Range("K136").Value = 14445
$ws.Range("L136").Value = 40900.5
$ws.Range("M136").Value = -11895
$ws.Range("N136").Value = -46000.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 114.86667
$ws.Range("I2").Value = 109.35714
$ws.Range("J2").Value = 119.6875
$ws.Range("K2").Value = 656.14284
$ws.Range("L2").Value = 718.125
$ws.Range("M2").Value = -543.14284
$ws.Range("N2").Value = -944.125
$ws.Range("H18").Value = 397.6154
$ws.Range("I18").Value = 178.25
$ws.Range("J18").Value = 3030
$ws.Range("K18").Value = 534.75
$ws.Range("L18").Value = 9090
$ws.Range("M18").Value = -365.75
$ws.Range("N18").Value = -9428
$ws.Range("H34").Value = 3527.7778
$ws.Range("I34").Value = 266.66666
$ws.Range("J34").Value = 5158.3335
$ws.Range("K34").Value = 799.9999799999999
$ws.Range("L34").Value = 15475.0005
$ws.Range("M34").Value = -715.9999799999999
$ws.Range("N34").Value = -15643.0005
$ws.Range("H36").Value = 1750
$ws.Range("I36").Value = 900
$ws.Range("J36").Value = 2033.3334
$ws.Range("K36").Value = 2700
$ws.Range("L36").Value = 6100.0002
$ws.Range("M36").Value = -2531
$ws.Range("N36").Value = -6438.0002
$ws.Range("H40").Value = 1147
$ws.Range("I40").Value = 1265
$ws.Range("J40").Value = 85
$ws.Range("K40").Value = 5060
$ws.Range("L40").Value = 340
$ws.Range("M40").Value = -4991
$ws.Range("N40").Value = -478
$ws.Range("H68").Value = 3554.423
$ws.Range("I68").Value = 2842.75
$ws.Range("J68").Value = 3870.7222
$ws.Range("K68").Value = 8528.25
$ws.Range("L68").Value = 11612.1666
$ws.Range("M68").Value = -7717.25
$ws.Range("N68").Value = -13234.1666
$ws.Range("H70").Value = 16816.666
$ws.Range("H71").Value = 3554.423
$ws.Range("I71").Value = 2842.75
$ws.Range("J71").Value = 3870.7222
$ws.Range("K71").Value = 25584.75
$ws.Range("L71").Value = 34836.49980000001
$ws.Range("M71").Value = -21528.75
$ws.Range("N71").Value = -42948.49980000001
$ws.Range("H73").Value = 16816.666
$ws.Range("H75").Value = 9085
$ws.Range("I75").Value = 874.3333
$ws.Range("K75").Value = 2622.9999
$ws.Range("M75").Value = -1624.9999
$ws.Range("H78").Value = 9085
$ws.Range("I78").Value = 874.3333
$ws.Range("K78").Value = 7868.9997
$ws.Range("M78").Value = -2876.9997
$ws.Range("H92").Value = 75.6
$ws.Range("I92").Value = 76.666664
$ws.Range("K92").Value = 229.999992
$ws.Range("M92").Value = 1018.000008
$ws.Range("H107").Value = 1769.9333
$ws.Range("I107").Value = 1625.8572
$ws.Range("J107").Value = 1896
$ws.Range("K107").Value = 4877.571599999999
$ws.Range("L107").Value = 5688
$ws.Range("M107").Value = -2957.571599999999
$ws.Range("N107").Value = -9528
$ws.Range("H112").Value = 13061
$ws.Range("I112").Value = 9374.5
$ws.Range("J112").Value = 17976.334
$ws.Range("K112").Value = 28123.5
$ws.Range("L112").Value = 53929.00199999999
$ws.Range("M112").Value = -27015.5
$ws.Range("N112").Value = -56145.00199999999
$ws.Range("H118").Value = 3694.75
$ws.Range("I118").Value = 3694.75
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 11084.25
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -9841.25
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 927.4545
$ws.Range("J122").Value = 921.8333
$ws.Range("L122").Value = 8296.4997
$ws.Range("N122").Value = -13196.4997
$ws.Range("H131").Value = 3032219.2
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1228.8667
$ws.Range("I2").Value = 1202.3572
$ws.Range("J2").Value = 1600
$ws.Range("K2").Value = 1202.3572
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = -1089.3572
$ws.Range("N2").Value = -1826
$ws.Range("H31").Value = 12896.667
$ws.Range("I31").Value = 8595
$ws.Range("K31").Value = 8595
$ws.Range("M31").Value = -8303
$ws.Range("H37").Value = 12896.667
$ws.Range("I37").Value = 8595
$ws.Range("K37").Value = 8595
$ws.Range("M37").Value = -8318
$ws.Range("H80").Value = 10382.385
$ws.Range("I80").Value = 8025.2856
$ws.Range("J80").Value = 13132.333
$ws.Range("K80").Value = 8025.2856
$ws.Range("L80").Value = 13132.333
$ws.Range("M80").Value = -7027.2856
$ws.Range("N80").Value = -15128.333
$ws.Range("H83").Value = 10382.385
$ws.Range("I83").Value = 8025.2856
$ws.Range("J83").Value = 13132.333
$ws.Range("K83").Value = 40126.428
$ws.Range("L83").Value = 65661.66500000001
$ws.Range("M83").Value = -35134.428
$ws.Range("N83").Value = -75645.66500000001
$ws.Range("H97").Value = 442.85715
$ws.Range("I97").Value = 442.85715
$ws.Range("K97").Value = 442.85715
$ws.Range("M97").Value = 53.14285000000001
$ws.Range("H102").Value = 2187.6829
$ws.Range("I102").Value = 1394.0286
$ws.Range("K102").Value = 1394.0286
$ws.Range("M102").Value = 227.9713999999999
$ws.Range("H122").Value = 2975.8975
$ws.Range("I122").Value = 2693.9644
$ws.Range("K122").Value = 8081.8932
$ws.Range("M122").Value = -5631.8932
$ws.Range("H132").Value = 2318.234
$ws.Range("I132").Value = 1987.6167
$ws.Range("J132").Value = 3485.1177
$ws.Range("K132").Value = 5962.8501
$ws.Range("L132").Value = 10455.3531
$ws.Range("M132").Value = -3432.8501
$ws.Range("N132").Value = -15515.3531

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 25000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 25000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 25000
$ws.Range("N18").Value = -25344
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 6000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10452
$ws.Range("H39").Value = 52000
$ws.Range("J39").Value = 52000
$ws.Range("L39").Value = 52000
$ws.Range("N39").Value = -52920
$ws.Range("H40").Value = 3341.2666
$ws.Range("I40").Value = 2176.6667
$ws.Range("K40").Value = 2176.6667
$ws.Range("M40").Value = -2040.6667
$ws.Range("H44").Value = 19750
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 19750
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 19750
$ws.Range("N44").Value = -20662
$ws.Range("M44").ClearContents()
$ws.Range("H46").Value = 4097.8335
$ws.Range("I46").Value = 3216
$ws.Range("J46").Value = 4538.75
$ws.Range("K46").Value = 3216
$ws.Range("L46").Value = 4538.75
$ws.Range("M46").Value = -3028
$ws.Range("N46").Value = -4914.75
$ws.Range("H61").Value = 2080.8647
$ws.Range("I61").Value = 1453.3334
$ws.Range("K61").Value = 1453.3334
$ws.Range("M61").Value = -1251.3334
$ws.Range("H93").Value = 1674.4667
$ws.Range("I93").Value = 1674.4667
$ws.Range("K93").Value = 1674.4667
$ws.Range("M93").Value = -426.4666999999999
$ws.Range("H111").Value = 116633
$ws.Range("J111").Value = 116633
$ws.Range("L111").Value = 116633
$ws.Range("N111").Value = -124813
$ws.Range("H113").Value = 2080.8647
$ws.Range("I113").Value = 1453.3334
$ws.Range("K113").Value = 1453.3334
$ws.Range("M113").Value = 716.6666
$ws.Range("H132").Value = 12363.23
$ws.Range("I132").Value = 6444
$ws.Range("J132").Value = 13439.454
$ws.Range("K132").Value = 19332
$ws.Range("L132").Value = 40318.362
$ws.Range("M132").Value = -16802
$ws.Range("N132").Value = -45378.362
$ws.Range("H136").Value = 7684.8213
$ws.Range("I136").Value = 7135.9165
$ws.Range("J136").Value = 8096.5
$ws.Range("K136").Value = 21407.7495
$ws.Range("L136").Value = 24289.5
$ws.Range("M136").Value = -18857.7495
$ws.Range("N136").Value = -29389.5
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 8000
$ws.Range("I32").Value = 8000
$ws.Range("K32").Value = 8000
$ws.Range("M32").Value = -7683
$ws.Range("H81").Value = 2009.2
$ws.Range("I81").Value = 2009.2
$ws.Range("K81").Value = 4018.4
$ws.Range("M81").Value = -2957.4
$ws.Range("H84").Value = 2009.2
$ws.Range("I84").Value = 2009.2
$ws.Range("K84").Value = 20092
$ws.Range("M84").Value = -14788
$ws.Range("H100").Value = 649.9091
$ws.Range("I100").Value = 537.5
$ws.Range("K100").Value = 1075
$ws.Range("M100").Value = -534
$ws.Range("H113").Value = 1633.238
$ws.Range("I113").Value = 574.3333
$ws.Range("J113").Value = 4280.5
$ws.Range("K113").Value = 1722.9999
$ws.Range("L113").Value = 12841.5
$ws.Range("M113").Value = 447.0001
$ws.Range("N113").Value = -17181.5
$ws.Range("H122").Value = 4231.3706
$ws.Range("I122").Value = 3845.5217
$ws.Range("K122").Value = 11536.5651
$ws.Range("M122").Value = -9086.5651
$ws.Range("H135").Value = 173000
$ws.Range("J135").Value = 173000
$ws.Range("L135").Value = 173000
$ws.Range("N135").Value = -183140
$ws.Range("H136").Value = 6429.879
$ws.Range("I136").Value = 4972.4546
$ws.Range("K136").Value = 14917.3638
$ws.Range("M136").Value = -12367.3638

